$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp text (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 3 de Julio de 2020 a las 17:15"

# --- Update case statistics for countries whose numbers changed but kept their row ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 2841055
$ws.Range("C4").Value = 3866
$ws.Range("D4").Value = 1191892
$ws.Range("E4").Value = 1517597
$ws.Range("G4").Value = 81
$ws.Range("H4").Value = 131566

# Row 7: India
$ws.Range("B7").Value = 633381
$ws.Range("C7").Value = 6213
$ws.Range("D7").Value = 383936
$ws.Range("E7").Value = 231125
$ws.Range("G7").Value = 95
$ws.Range("H7").Value = 18320

# Row 29: Bielorrusia
$ws.Range("B29").Value = 62997
$ws.Range("C29").Value = 299
$ws.Range("D29").Value = 49909
$ws.Range("E29").Value = 12676
$ws.Range("G29").Value = 7
$ws.Range("H29").Value = 412

# Row 50: Barein
$ws.Range("E50").Value = 5159
$ws.Range("G50").Value = 1
$ws.Range("H50").Value = 95

# Row 61: Moldavia
$ws.Range("D61").Value = 10093
$ws.Range("E61").Value = 6494
$ws.Range("G61").Value = 3
$ws.Range("H61").Value = 563

# Rows 81/82: Kirguistan and Republica de Macedonia swap order because
# Macedonia's updated total (6787) overtakes Kirguistan's (6767).
$ws.Range("A81").Value = "Republica de Macedonia"
$ws.Range("B81").Value = 6787
$ws.Range("C81").Value = 162
$ws.Range("D81").Value = 2876
$ws.Range("E81").Value = 3583
$ws.Range("G81").Value = 7
$ws.Range("H81").Value = 328

$ws.Range("A82").Value = "Kirguistan"
$ws.Range("B82").Value = 6767
$ws.Range("C82").Value = 506
$ws.Range("D82").Value = 2655
$ws.Range("E82").Value = 4036
$ws.Range("G82").Value = 10
$ws.Range("H82").Value = 76

# Row 93: Guayana Francesa
$ws.Range("B93").Value = 4558
$ws.Range("C93").Value = 114
$ws.Range("D93").Value = 1777
$ws.Range("E93").Value = 2765

# Row 105: Mayotte
$ws.Range("B105").Value = 2661
$ws.Range("C105").Value = 11
$ws.Range("D105").Value = 2375
$ws.Range("E105").Value = 251

# Row 115: Islandia
$ws.Range("B115").Value = 1855
$ws.Range("C115").Value = 5
$ws.Range("E115").Value = 17

# Row 164: Birmania
$ws.Range("B164").Value = 306
$ws.Range("C164").Value = 2
$ws.Range("D164").Value = 237
$ws.Range("E164").Value = 63

# Rows 205/206: Fiyi and Dominica swap order (tied totals, reordered)
$ws.Range("A205").Value = "Dominica"
$ws.Range("A206").Value = "Fiyi"

# Rows 209/210: Groenlandia and Islas Malvinas swap order (tied totals, reordered)
$ws.Range("A209").Value = "Islas Malvinas"
$ws.Range("A210").Value = "Groenlandia"
